{"js": "// Remove the empty paragraph and the long \"Creative Commons ... no es una\n// parte en sus licencias p\u00fablicas ...\" disclaimer paragraph that sit right\n// after the second horizontal-rule picture, just before the \"Comunes\n// Creativos puede ser contactada en ...\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the unique paragraph whose text starts the long CC disclaimer.\nconst marker = \"Creative Commons no es una parte en sus licencias p\u00fablicas\";\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) === 0) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not find the target 'Creative Commons' disclaimer paragraph.\");\n}\n\nconst targetParagraph = paragraphs.items[targetIndex];\nconst previousParagraph = paragraphs.items[targetIndex - 1];\n\n// The paragraph immediately before the target one is the empty paragraph\n// that follows the horizontal rule picture; confirm it is blank before\n// removing it so we don't delete the wrong paragraph if structure shifts.\npreviousParagraph.load(\"text\");\nawait context.sync();\n\nif (previousParagraph.text.trim() === \"\") {\n  previousParagraph.delete();\n}\n\ntargetParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Remove the empty paragraph and the long \"Creative Commons ... no es una\n# parte en sus licencias p\u00fablicas ...\" disclaimer paragraph that sit right\n# after the second horizontal-rule picture, just before the \"Comunes\n# Creativos puede ser contactada en ...\" paragraph.\n\n$d = $word.ActiveDocument\n\n$marker = \"Creative Commons no es una parte en sus licencias p\u00fablicas\"\n\n$paragraphs = $d.Paragraphs\n$targetIndex = -1\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $paraText = $paragraphs.Item($i).Range.Text\n    if ($paraText.IndexOf($marker) -eq 0) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the target 'Creative Commons' disclaimer paragraph.\"\n}\n\n# Delete the long disclaimer paragraph first (this merges it away cleanly).\n$targetParagraph = $d.Paragraphs.Item($targetIndex)\n$targetParagraph.Range.Delete()\n\n# The paragraph that is now at (targetIndex - 1) is the empty paragraph that\n# used to sit between the horizontal-rule picture and the disclaimer; remove\n# it too, but only if it is in fact blank so we don't delete real content.\n$previousParagraph = $d.Paragraphs.Item($targetIndex - 1)\nif ($previousParagraph.Range.Text.Trim() -eq \"\") {\n    $previousParagraph.Range.Delete()\n}\n"}
